$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 4912.87
$ws.Range("C2").Value = 11208.15
$ws.Range("D2").Value = 2031
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 23697.44
$ws.Range("G2").Value = 10948.72
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 0
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 0
$ws.Range("AG2").Value = 52798.18

# Row 3
$ws.Range("B3").Value = 4803
$ws.Range("C3").Value = 4651.9
$ws.Range("D3").Value = 19192
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 5661
$ws.Range("G3").Value = 5242.9
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = 0
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = 0
$ws.Range("AC3").Value = 0
$ws.Range("AD3").Value = 0
$ws.Range("AG3").Value = 39550.8

# Row 4
$ws.Range("B4").Value = 5020
$ws.Range("C4").Value = 2670
$ws.Range("D4").Value = 4432
$ws.Range("E4").Value = 2250
$ws.Range("F4").Value = 3344.5
$ws.Range("G4").Value = 679
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 0
$ws.Range("Y4").Value = 0
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0
$ws.Range("AC4").Value = 0
$ws.Range("AD4").Value = 0
$ws.Range("AG4").Value = 18395.5

# Row 5
$ws.Range("B5").Value = 2321
$ws.Range("C5").Value = 4256.58
$ws.Range("D5").Value = 3600
$ws.Range("E5").Value = 2120
$ws.Range("F5").Value = 2699.7
$ws.Range("G5").Value = 864.8
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0
$ws.Range("Y5").Value = 0
$ws.Range("Z5").Value = 0
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0
$ws.Range("AC5").Value = 0
$ws.Range("AD5").Value = 0
$ws.Range("AG5").Value = 15862.08

# Row 6
$ws.Range("B6").Value = 17056.87
$ws.Range("C6").Value = 22786.63
$ws.Range("D6").Value = 29255
$ws.Range("E6").Value = 4370
$ws.Range("F6").Value = 35402.64
$ws.Range("G6").Value = 17735.42
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 0
$ws.Range("Y6").Value = 0
$ws.Range("Z6").Value = 0
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0
$ws.Range("AC6").Value = 0
$ws.Range("AD6").Value = 0
$ws.Range("AG6").Value = 126606.56
